$d = $word.ActiveDocument

# Locate the unique sentence fragment "weblogs_flat WHERE LIMIT 10;" so we don't
# accidentally touch the other (legitimate) WHERE clause later in the document.
$sentence = $d.Content
$foundSentence = $sentence.Find.Execute("weblogs_flat WHERE LIMIT 10;", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundSentence) {
    throw "Could not locate the target sentence fragment"
}

# Within that fragment, find the exact "WHERE " text (including the trailing
# space) that needs to be removed as part of the typo fix.
$whereRange = $d.Range($sentence.Start, $sentence.End)
$foundWhere = $whereRange.Find.Execute("WHERE ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundWhere) {
    throw "Could not locate 'WHERE ' within the target sentence"
}
$whereStart = $whereRange.Start
$whereEnd = $whereRange.End

# Word keeps the "_GoBack" bookmark pinned to the location of the most recent
# edit. Re-anchor it (collapsed) at the point where "WHERE " is about to be
# removed *before* performing the deletion, mirroring how Word itself updates
# the bookmark as the user edits text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$goBackAnchor = $d.Range($whereStart, $whereStart)
$goBackAnchor.Bookmarks.Add("_GoBack")

# Now remove the stray "WHERE " text, leaving "LIMIT 10;" directly after the
# preceding space.
$d.Range($whereStart, $whereEnd).Delete()
